# Auto-generated edit script applying the Garuda_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H24").Value = 500
$ws.Range("J24").Value = 500
$ws.Range("L24").Value = 1500
$ws.Range("N24").Value = -1840
$ws.Range("H132").Value = 8405738
$ws.Range("I132").Value = 8405738
$ws.Range("K132").Value = 25217214
$ws.Range("M132").Value = -25214684
$ws.Range("H135").Value = 1037.4546
$ws.Range("I135").Value = 1038.625
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 9347.625
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -6812.625
$ws.Range("N135").Value = -14070

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2719.2
$ws.Range("I2").Value = 1299.5
$ws.Range("J2").Value = 3665.6667
$ws.Range("K2").Value = 1299.5
$ws.Range("L2").Value = 3665.6667
$ws.Range("M2").Value = -1186.5
$ws.Range("N2").Value = -3891.6667
$ws.Range("H45").Value = 856.375
$ws.Range("I45").Value = 783.82355
$ws.Range("J45").Value = 938.6
$ws.Range("K45").Value = 783.82355
$ws.Range("L45").Value = 938.6
$ws.Range("M45").Value = -406.82355
$ws.Range("N45").Value = -1692.6
$ws.Range("H61").Value = 1732.4783
$ws.Range("I61").Value = 935.5333000000001
$ws.Range("J61").Value = 3226.75
$ws.Range("K61").Value = 935.5333000000001
$ws.Range("L61").Value = 3226.75
$ws.Range("M61").Value = -723.5333000000001
$ws.Range("N61").Value = -3650.75
$ws.Range("H88").Value = 775616.9
$ws.Range("I88").Value = 1255001.5
$ws.Range("J88").Value = 8601.4
$ws.Range("K88").Value = 1255001.5
$ws.Range("L88").Value = 8601.4
$ws.Range("M88").Value = -1254595.5
$ws.Range("N88").Value = -9413.4
$ws.Range("H91").Value = 775616.9
$ws.Range("I91").Value = 1255001.5
$ws.Range("J91").Value = 8601.4
$ws.Range("K91").Value = 1255001.5
$ws.Range("L91").Value = 8601.4
$ws.Range("M91").Value = -1253597.5
$ws.Range("N91").Value = -11409.4
$ws.Range("H116").Value = 2719.2
$ws.Range("I116").Value = 1299.5
$ws.Range("J116").Value = 3665.6667
$ws.Range("K116").Value = 1299.5
$ws.Range("L116").Value = 3665.6667
$ws.Range("M116").Value = 994.5
$ws.Range("N116").Value = -8253.6667
$ws.Range("H122").Value = 1006
$ws.Range("I122").Value = 1006
$ws.Range("K122").Value = 3018
$ws.Range("M122").Value = -568
$ws.Range("H136").Value = 1732.4783
$ws.Range("I136").Value = 935.5333000000001
$ws.Range("J136").Value = 3226.75
$ws.Range("K136").Value = 2806.5999
$ws.Range("L136").Value = 9680.25
$ws.Range("M136").Value = -256.5999000000002
$ws.Range("N136").Value = -14780.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2719.2
$ws.Range("I3").Value = 1299.5
$ws.Range("J3").Value = 3665.6667
$ws.Range("K3").Value = 1299.5
$ws.Range("L3").Value = 3665.6667
$ws.Range("M3").Value = -1185.5
$ws.Range("N3").Value = -3893.6667
$ws.Range("H20").Value = 3554.2258
$ws.Range("I20").Value = 2173.3125
$ws.Range("J20").Value = 5027.2
$ws.Range("K20").Value = 2173.3125
$ws.Range("L20").Value = 5027.2
$ws.Range("M20").Value = -1926.3125
$ws.Range("N20").Value = -5521.2
$ws.Range("H80").Value = 628
$ws.Range("J80").Value = 520.75
$ws.Range("L80").Value = 520.75
$ws.Range("N80").Value = -2516.75
$ws.Range("H83").Value = 628
$ws.Range("J83").Value = 520.75
$ws.Range("L83").Value = 2603.75
$ws.Range("N83").Value = -12587.75
$ws.Range("H86").Value = 2243
$ws.Range("I86").Value = 2298.3333
$ws.Range("J86").Value = 2201.5
$ws.Range("K86").Value = 2298.3333
$ws.Range("L86").Value = 2201.5
$ws.Range("M86").Value = -1175.3333
$ws.Range("N86").Value = -4447.5
$ws.Range("H89").Value = 2243
$ws.Range("I89").Value = 2298.3333
$ws.Range("J89").Value = 2201.5
$ws.Range("K89").Value = 11491.6665
$ws.Range("L89").Value = 11007.5
$ws.Range("M89").Value = -5875.666499999999
$ws.Range("N89").Value = -22239.5
$ws.Range("H99").Value = 950.6111
$ws.Range("I99").Value = 920
$ws.Range("J99").Value = 988.875
$ws.Range("K99").Value = 920
$ws.Range("L99").Value = 988.875
$ws.Range("M99").Value = 578
$ws.Range("N99").Value = -3984.875

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8698376
$ws.Range("I31").Value = 2734.25
$ws.Range("J31").Value = 18184532
$ws.Range("K31").Value = 2734.25
$ws.Range("L31").Value = 18184532
$ws.Range("M31").Value = -2439.25
$ws.Range("N31").Value = -18185122
$ws.Range("H34").Value = 8698376
$ws.Range("I34").Value = 2734.25
$ws.Range("J34").Value = 18184532
$ws.Range("K34").Value = 2734.25
$ws.Range("L34").Value = 18184532
$ws.Range("M34").Value = -2532.25
$ws.Range("N34").Value = -18184936
$ws.Range("H58").Value = 890.62964
$ws.Range("I58").Value = 961.17645
$ws.Range("J58").Value = 770.7
$ws.Range("K58").Value = 961.17645
$ws.Range("L58").Value = 770.7
$ws.Range("M58").Value = -758.17645
$ws.Range("N58").Value = -1176.7
$ws.Range("H62").Value = 2664.6667
$ws.Range("I62").Value = 2499
$ws.Range("J62").Value = 2996
$ws.Range("K62").Value = 2499
$ws.Range("L62").Value = 2996
$ws.Range("M62").Value = -1875
$ws.Range("N62").Value = -4244
$ws.Range("H65").Value = 2664.6667
$ws.Range("I65").Value = 2499
$ws.Range("J65").Value = 2996
$ws.Range("K65").Value = 12495
$ws.Range("L65").Value = 14980
$ws.Range("M65").Value = -9375
$ws.Range("N65").Value = -21220
$ws.Range("H107").Value = 486.34885
$ws.Range("I107").Value = 461.22726
$ws.Range("J107").Value = 512.6667
$ws.Range("K107").Value = 461.22726
$ws.Range("L107").Value = 512.6667
$ws.Range("M107").Value = 1458.77274
$ws.Range("N107").Value = -4352.6667
$ws.Range("H134").Value = 1282.1428
$ws.Range("I134").Value = 1282.1428
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3846.4284
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1311.4284
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 890.62964
$ws.Range("I136").Value = 961.17645
$ws.Range("J136").Value = 770.7
$ws.Range("K136").Value = 2883.52935
$ws.Range("L136").Value = 2312.1
$ws.Range("M136").Value = -333.5293500000002
$ws.Range("N136").Value = -7412.1

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 554.03845
$ws.Range("J5").Value = 738.5625
$ws.Range("L5").Value = 2215.6875
$ws.Range("N5").Value = -2439.6875
$ws.Range("H131").Value = 2658111.5
$ws.Range("J131").Value = 11905212
$ws.Range("L131").Value = 35715636
$ws.Range("N131").Value = -35725716
$ws.Range("H135").Value = 554.03845
$ws.Range("J135").Value = 738.5625
$ws.Range("L135").Value = 6647.0625
$ws.Range("N135").Value = -11717.0625

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 31879462
$ws.Range("I70").Value = 36432956
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 36432956
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -36432686
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 31879462
$ws.Range("I73").Value = 36432956
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 36432956
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -36432020
$ws.Range("N73").Value = -6872
$ws.Range("H102").Value = 1542
$ws.Range("I102").Value = 1480.8182
$ws.Range("J102").Value = 1766.3334
$ws.Range("K102").Value = 1480.8182
$ws.Range("L102").Value = 1766.3334
$ws.Range("M102").Value = 141.1818000000001
$ws.Range("N102").Value = -5010.3334
$ws.Range("H122").Value = 2616.8
$ws.Range("I122").Value = 2472.5334
$ws.Range("K122").Value = 7417.600199999999
$ws.Range("M122").Value = -4967.600199999999
$ws.Range("H126").Value = 2510.1904
$ws.Range("I126").Value = 2690
$ws.Range("J126").Value = 2346.7273
$ws.Range("K126").Value = 8070
$ws.Range("L126").Value = 7040.1819
$ws.Range("M126").Value = -5600
$ws.Range("N126").Value = -11980.1819

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4863.84
$ws.Range("I16").Value = 890.7143
$ws.Range("J16").Value = 9920.546
$ws.Range("K16").Value = 890.7143
$ws.Range("L16").Value = 9920.546
$ws.Range("M16").Value = -720.7143
$ws.Range("N16").Value = -10260.546
$ws.Range("H133").Value = 14989.765
$ws.Range("J133").Value = 14989.765
$ws.Range("L133").Value = 14989.765
$ws.Range("N133").Value = -20049.765

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 20014
$ws.Range("J26").Value = 20014
$ws.Range("L26").Value = 20014
$ws.Range("N26").Value = -20600
$ws.Range("H136").Value = 5644.077
$ws.Range("I136").Value = 5644.077
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 16932.231
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -14382.231
$ws.Range("N136").ClearContents()
